# Add sheets via upload: recreates the sheet layout, shared-string data,
# and active-tab selection described by the target workbook.

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd {
    param($wb, $name)
    $count = $wb.Worksheets.Count
    $last = $wb.Worksheets.Item($count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $newSheet.Name = $name
    return $newSheet
}

# Helper: write a numeric-looking value as TEXT (not auto-coerced to a
# number) without leaving a stray NumberFormat style behind on the cell -
# format the cell as text, assign the value, then paste-special the
# (default) number format back in from a neighboring plain-text cell.
function Set-TextValue {
    param($range, $formatDonor, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $formatDonor.Copy()
    $range.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

$s3  = Add-SheetAtEnd $wb "1"
$s4  = Add-SheetAtEnd $wb "061-FEB-2018"
$s5  = Add-SheetAtEnd $wb "06-FEB-2018"
$s6  = Add-SheetAtEnd $wb "06-FEB-20181"
$s7  = Add-SheetAtEnd $wb "07-FEB-2018"
$s8  = Add-SheetAtEnd $wb "07-FEB-20181"
$s9  = Add-SheetAtEnd $wb "07-FEB-20182"
$s10 = Add-SheetAtEnd $wb "07-FEB-20183"

# Sheet "1": a lone "test" string in B2
$s3.Range("B2").Value = "test"

# Sheet "061-FEB-2018": header row + one data row
$s4.Range("B1").Value = "Transaction Reference"
$s4.Range("C1").Value = "Service Agreement"
$s4.Range("D1").Value = "MAC"
$s4.Range("B2").Value = "SPark_06FEB20183"
$s4.Range("D2").Value = "AA:BB:CC:2B:F4:FB"
Set-TextValue $s4.Range("C2") $s4.Range("D2") "402906998"

# Sheet "06-FEB-2018": header row + one data row
$s5.Range("B1").Value = "Transaction Reference"
$s5.Range("C1").Value = "Service Agreement"
$s5.Range("D1").Value = "MAC"
$s5.Range("B2").Value = "SPark_06FEB20181"
$s5.Range("D2").Value = "AA:BB:CC:3B:09:51"
Set-TextValue $s5.Range("C2") $s5.Range("D2") "402907168"

# Sheet "06-FEB-20181": left empty (no data)

# Sheet "07-FEB-2018": header row + one data row
$s7.Range("B1").Value = "Transaction Reference"
$s7.Range("C1").Value = "Service Agreement"
$s7.Range("D1").Value = "MAC"
$s7.Range("B2").Value = "SPark_07FEB20183"
$s7.Range("D2").Value = "AA:BB:CC:01:36:C2"
Set-TextValue $s7.Range("C2") $s7.Range("D2") "402907738"

# Sheets "07-FEB-20181", "07-FEB-20182", "07-FEB-20183": left empty

# The newly-added "061-FEB-2018" sheet (4th tab, 0-indexed activeTab=3)
# ends up the selected/active tab, with F32 as the last selected cell.
$s4.Select()
$s4.Range("F32").Select() | Out-Null
